# Apply updated cryptocurrency market data values to worksheet 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price cells contain values that look like plain decimal numbers
# (e.g. "214.44"). Force those specific cells to Text format first so
# Excel stores them as strings (matching the source data) instead of
# silently converting them to numeric values.
$textCells = @("D5", "D10", "D15", "D18", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D33", "D34", "D37", "D41", "D43", "D46", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.968.54'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.635.50'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '214.44'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("E9").Value = '  -2.99%  '
$ws.Range("D10").Value = '18.50'
$ws.Range("E10").Value = '  -6.09%  '
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '1.863.11'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.637.00'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("D15").Value = '0.533'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '25.993.40'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0₃0744'
$ws.Range("E17").Value = '  -2.95%  '
$ws.Range("D18").Value = '61.69'
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = '190.92'
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").Value = '9.62'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("D23").Value = '6.14'
$ws.Range("E23").Value = '  -2.00%  '
$ws.Range("D24").Value = '0.133'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = '143.10'
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = '1.77'
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").Value = '6.81'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  -3.46%  '
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("D33").Value = '3.15'
$ws.Range("E33").Value = '  -4.37%  '
$ws.Range("D34").Value = '2.43'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("D36").Value = '1.137.82'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").Value = '0.864'
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("E39").Value = '  -3.82%  '
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("D41").Value = '98.62'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.24'
$ws.Range("E43").Value = '  -5.09%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.772.62'
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").Value = '55.26'
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").Value = '  +1.87%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = '7.54'
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  -0.02%  '
